# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2036.909
$ws.Range("I19").Value = 1892.4286
$ws.Range("K19").Value = 1892.4286
$ws.Range("M19").Value = -1717.4286

$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2825

$ws.Range("H98").Value = 1569.8572
$ws.Range("I98").Value = 1330.5
$ws.Range("K98").Value = 1330.5
$ws.Range("M98").Value = 167.5

$ws.Range("H122").Value = 1569.8572
$ws.Range("I122").Value = 1330.5
$ws.Range("K122").Value = 3991.5
$ws.Range("M122").Value = -1541.5

$ws.Range("H129").Value = 1180.5264
$ws.Range("I129").Value = 691.2
$ws.Range("K129").Value = 2073.6
$ws.Range("M129").Value = 2926.4

$ws.Range("H138").Value = 2368.4644
$ws.Range("I138").Value = 1332.7
$ws.Range("J138").Value = 2943.889
$ws.Range("K138").Value = 3998.1
$ws.Range("L138").Value = 8831.667000000001
$ws.Range("M138").Value = 1141.9
$ws.Range("N138").Value = -19111.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 252.5
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 500
$ws.Range("M3").Value = 110
$ws.Range("N3").Value = -730

$ws.Range("H61").Value = 1986.7778
$ws.Range("I61").Value = 1986.7778
$ws.Range("K61").Value = 1986.7778
$ws.Range("M61").Value = -1774.7778

$ws.Range("H74").Value = 2277
$ws.Range("I74").Value = 1582.1666
$ws.Range("K74").Value = 1582.1666
$ws.Range("M74").Value = -708.1666

$ws.Range("H77").Value = 2277
$ws.Range("I77").Value = 1582.1666
$ws.Range("K77").Value = 7910.833000000001
$ws.Range("M77").Value = -3542.833000000001

$ws.Range("H110").Value = 1184
$ws.Range("I110").Value = 1184
$ws.Range("K110").Value = 1184
$ws.Range("M110").Value = 861

$ws.Range("H136").Value = 1986.7778
$ws.Range("I136").Value = 1986.7778
$ws.Range("K136").Value = 5960.3334
$ws.Range("M136").Value = -3410.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1646.1904
$ws.Range("I94").Value = 1620.5883
$ws.Range("J94").Value = 1755
$ws.Range("K94").Value = 1620.5883
$ws.Range("L94").Value = 1755
$ws.Range("M94").Value = -1169.5883
$ws.Range("N94").Value = -2657

$ws.Range("H105").Value = 3024.5
$ws.Range("I105").Value = 2826.818
$ws.Range("K105").Value = 2826.818
$ws.Range("M105").Value = -1079.818

$ws.Range("H107").Value = 998
$ws.Range("I107").Value = 920.9231
$ws.Range("K107").Value = 920.9231
$ws.Range("M107").Value = 999.0769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.66667
$ws.Range("I7").Value = 83.14286
$ws.Range("J7").Value = 238
$ws.Range("K7").Value = 83.14286
$ws.Range("L7").Value = 238
$ws.Range("M7").Value = 29.85714
$ws.Range("N7").Value = -464

$ws.Range("H10").Value = 186.25
$ws.Range("I10").Value = 186.25
$ws.Range("K10").Value = 186.25
$ws.Range("M10").Value = -47.25

$ws.Range("H16").Value = 996.5454999999999
$ws.Range("I16").Value = 1096.3334
$ws.Range("J16").Value = 547.5
$ws.Range("K16").Value = 1096.3334
$ws.Range("L16").Value = 547.5
$ws.Range("M16").Value = -809.3334
$ws.Range("N16").Value = -1121.5

$ws.Range("H113").Value = 996.5454999999999
$ws.Range("I113").Value = 1096.3334
$ws.Range("J113").Value = 547.5
$ws.Range("K113").Value = 1096.3334
$ws.Range("L113").Value = 547.5
$ws.Range("M113").Value = 1073.6666
$ws.Range("N113").Value = -4887.5

$ws.Range("H132").Value = 1302.7059
$ws.Range("I132").Value = 1031.9286
$ws.Range("J132").Value = 2566.3333
$ws.Range("K132").Value = 3095.7858
$ws.Range("L132").Value = 7698.999899999999
$ws.Range("M132").Value = -565.7857999999997
$ws.Range("N132").Value = -12758.9999

$ws.Range("H134").Value = 2353.8948
$ws.Range("I134").Value = 1920.8667
$ws.Range("K134").Value = 5762.6001
$ws.Range("M134").Value = -3227.6001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1351
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1351
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4053
$ws.Range("N22").Value = -4391
$ws.Range("M22").ClearContents()

$ws.Range("H23").Value = 285
$ws.Range("I23").Value = 20
$ws.Range("K23").Value = 60
$ws.Range("M23").Value = 175

$ws.Range("H27").Value = 1351
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1351
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4053
$ws.Range("N27").Value = -4257
$ws.Range("M27").ClearContents()

$ws.Range("H68").Value = 1528.5714
$ws.Range("I68").Value = 1275
$ws.Range("J68").Value = 1866.6666
$ws.Range("K68").Value = 3825
$ws.Range("L68").Value = 5599.9998
$ws.Range("M68").Value = -3014
$ws.Range("N68").Value = -7221.9998

$ws.Range("H71").Value = 1528.5714
$ws.Range("I71").Value = 1275
$ws.Range("J71").Value = 1866.6666
$ws.Range("K71").Value = 11475
$ws.Range("L71").Value = 16799.9994
$ws.Range("M71").Value = -7419
$ws.Range("N71").Value = -24911.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3988
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 3988
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 44063.5
$ws.Range("I7").Value = 52036.4
$ws.Range("J7").Value = 4199
$ws.Range("K7").Value = 52036.4
$ws.Range("L7").Value = 4199
$ws.Range("M7").Value = -51924.4
$ws.Range("N7").Value = -4423

$ws.Range("H22").Value = 1400
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1400
$ws.Range("N22").Value = -1990
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 1400
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1400
$ws.Range("N27").Value = -1614
$ws.Range("M27").ClearContents()

$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2312

$ws.Range("H126").Value = 44063.5
$ws.Range("I126").Value = 52036.4
$ws.Range("J126").Value = 4199
$ws.Range("K126").Value = 156109.2
$ws.Range("L126").Value = 12597
$ws.Range("M126").Value = -153639.2
$ws.Range("N126").Value = -17537

$ws.Range("H132").Value = 4700.647
$ws.Range("I132").Value = 1979.7778
$ws.Range("K132").Value = 5939.3334
$ws.Range("M132").Value = -3409.3334

$ws.Range("H136").Value = 2445.7
$ws.Range("I136").Value = 2250.375
$ws.Range("J136").Value = 3227
$ws.Range("K136").Value = 6751.125
$ws.Range("L136").Value = 9681
$ws.Range("M136").Value = -4201.125
$ws.Range("N136").Value = -14781

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4958.3335
$ws.Range("I62").Value = 3916.6667
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 3916.6667
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -3292.6667
$ws.Range("N62").Value = -7248

$ws.Range("H65").Value = 4958.3335
$ws.Range("I65").Value = 3916.6667
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 19583.3335
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -16463.3335
$ws.Range("N65").Value = -36240

$ws.Range("H132").Value = 1329.5625
$ws.Range("I132").Value = 1321.8462
$ws.Range("J132").Value = 1363
$ws.Range("K132").Value = 3965.5386
$ws.Range("L132").Value = 4089
$ws.Range("M132").Value = -1435.5386
$ws.Range("N132").Value = -9149
